# Practiced Cambridge 14 Test 4 -> log a new row (row 20) in the IELTS score
# table on Sheet1 with the results of that practice test.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 20 is currently an "empty placeholder" row inside Table1, styled with
# the light/no-border placeholder formats (s=10 / s=11) used for Listening
# and Reading calculated columns. Once a row gets real data those columns
# pick up the "filled" look (s=8 / s=9) that the rows above already use, so
# copy that formatting down from row 19 (the previous entry) before typing
# the new values in.
$ws.Range("F19").Copy() | Out-Null
$ws.Range("F20").PasteSpecial(-4122) | Out-Null
$ws.Range("H19").Copy() | Out-Null
$ws.Range("H20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Enter the Cambridge IELTS 14 Test 4 practice results.
$ws.Range("C20").Value2 = 45466
$ws.Range("D20").Value2 = "IELTS14_Test4"
$ws.Range("E20").Value2 = 32
$ws.Range("F20").Formula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20, MATCH(Table1[[#This Row],[Lis_Mark]], Sheet2!$D$5:$D$20, 1)),"No Grade")'
$ws.Range("G20").Value2 = 26
$ws.Range("H20").Formula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20, MATCH(Table1[[#This Row],[Read_Mark]], Sheet2!$D$5:$D$20, 1)),"No Grade")'
$ws.Range("I20").Value2 = 1.1
$ws.Range("J20").Value2 = 3
$ws.Range("K20").Formula = "=(F20+H20+I20+J20)/4"

# Leave the selection where the user ended up after filling in the row.
$ws.Range("I23").Select() | Out-Null
